# Applies the cryptos.xlsx price/volume/coin-order refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields (coin name / link) - no numeric ambiguity, assign directly.
$textUpdates = @{
    'B28' = 'Kaspa'
    'C28' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'B29' = 'Dai'
    'C29' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'B33' = 'VeChain'
    'C33' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'B34' = 'OKB'
    'C34' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'B35' = 'InjectiveProtocol'
    'C35' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value2 = $textUpdates[$ref]
}

# Price / Volume(1h) fields are stored as literal text (e.g. "51.020.74", "5.00", "  -1.52%  ").
# Force text format before writing so Excel does not reinterpret these as numbers/percentages
# and silently trim formatting (trailing zeros, thousand-dot grouping, padding spaces),
# then restore the default "Normal" style so no stray number-format style sticks to the cell.
$numericLookingUpdates = @{
    'D2' = '51.020.74'
    'E2' = '  -1.52%  '
    'D3' = '2.935.28'
    'E3' = '  -2.60%  '
    'E4' = '  +0.09%  '
    'D5' = '376.99'
    'E5' = '  -1.42%  '
    'D6' = '102.18'
    'E6' = '  -4.25%  '
    'E7' = '  -2.84%  '
    'E8' = '  +0.08%  '
    'D9' = '0.583'
    'E9' = '  -3.59%  '
    'D10' = '36.51'
    'E10' = '  -4.07%  '
    'D11' = '0.138'
    'E11' = '  -1.20%  '
    'D12' = '0.0834'
    'E12' = '  -1.89%  '
    'D13' = '3.391.91'
    'E13' = '  -2.40%  '
    'D14' = '17.93'
    'E14' = '  -5.36%  '
    'D15' = '7.33'
    'E15' = '  -3.30%  '
    'D16' = '2.955.28'
    'E16' = '  -1.74%  '
    'D17' = '0.971'
    'E17' = '  -1.06%  '
    'D18' = '50.930.99'
    'E18' = '  -1.64%  '
    'D19' = '3.16'
    'E19' = '  -7.94%  '
    'D20' = '7.12'
    'E20' = '  -5.09%  '
    'D21' = '12.49'
    'E21' = '  -5.19%  '
    'D22' = '0.0₃0947'
    'E22' = '  -2.22%  '
    'D23' = '68.06'
    'E23' = '  -1.50%  '
    'D24' = '261.42'
    'E24' = '  -1.34%  '
    'D25' = '2.86'
    'E25' = '  +1.47%  '
    'D26' = '8.18'
    'E26' = '  +7.75%  '
    'D27' = '7.58'
    'E27' = '  +3.75%  '
    'D28' = '0.168'
    'E28' = '  -3.06%  '
    'D29' = '1.00'
    'E29' = '  +0.05%  '
    'E30' = '  +5.40%  '
    'D31' = '25.56'
    'E31' = '  -2.89%  '
    'D32' = '9.79'
    'E32' = '  -2.31%  '
    'D33' = '0.0455'
    'E33' = '  +0.95%  '
    'D34' = '50.62'
    'E34' = '  -1.13%  '
    'D35' = '33.80'
    'E35' = '  -3.84%  '
    'E36' = '  -3.17%  '
    'E37' = '  +0.00%  '
    'D38' = '2.97'
    'E38' = '  -5.26%  '
    'E39' = '  -4.52%  '
    'E40' = '  -2.65%  '
    'D41' = '16.29'
    'E41' = '  -7.99%  '
    'D42' = '1.78'
    'E42' = '  -5.71%  '
    'D43' = '121.47'
    'E43' = '  -2.56%  '
    'D44' = '21.04'
    'E44' = '  -6.86%  '
    'E45' = '  -1.89%  '
    'E46' = '  -2.71%  '
    'D47' = '0.271'
    'E47' = '  -2.19%  '
    'D48' = '2.002.15'
    'E48' = '  -2.98%  '
    'D49' = '3.21'
    'E49' = '  -3.07%  '
    'D50' = '0.0343'
    'E50' = '  -3.99%  '
    'D51' = '5.00'
    'E51' = '  -4.67%  '
}
foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value2 = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}
